# edit.ps1
# Applies the weekly CompStat update (new crime data) described by the diff:
# - Report header text: Volume .. Number 47 -> 48; week dates 11/21/2022-11/27/2022 -> 11/28/2022-12/4/2022
# - Numeric crime-complaint figures for rows 14-30 (precincts / boroughs section)
# - Row 30 (Housing borough) gains real Week-to-Date figures for the "Prior" (D) and
#   "% Chg" (E) columns, switching those two cells from placeholder text to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text cells; re-set as plain text, the cell-level
#     font style already matches the per-run formatting so appearance is unchanged) ---
$ws.Range("A8").Value = "Volume 29   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/28/2022  Through  12/4/2022"

# --- Row 30 special case: D30/E30 move from placeholder text to real numbers ---
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 3
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E30").Value = -100

# --- Bulk numeric updates ---
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = -66.666666666666
$ws.Range("I14").Value = 116
$ws.Range("J14").Value = 137
$ws.Range("K14").Value = -15.328467153284
$ws.Range("L14").Value = 10.47619047619
$ws.Range("M14").Value = -4.132231404958
$ws.Range("N14").Value = -75.371549893842
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 9
$ws.Range("E15").Value = -11.111111111111
$ws.Range("F15").Value = 22
$ws.Range("H15").Value = -15.384615384615
$ws.Range("I15").Value = 354
$ws.Range("J15").Value = 342
$ws.Range("K15").Value = 3.508771929824
$ws.Range("L15").Value = 15.686274509803
$ws.Range("M15").Value = 27.797833935018
$ws.Range("N15").Value = -47.399702823179
$ws.Range("C16").Value = 83
$ws.Range("D16").Value = 111
$ws.Range("E16").Value = -25.225225225225
$ws.Range("F16").Value = 398
$ws.Range("G16").Value = 422
$ws.Range("H16").Value = -5.687203791469
$ws.Range("I16").Value = 4756
$ws.Range("J16").Value = 3679
$ws.Range("K16").Value = 29.274259309595
$ws.Range("L16").Value = 47.518610421836
$ws.Range("M16").Value = 10.888318955467
$ws.Range("N16").Value = -68.74548202668
$ws.Range("C17").Value = 118
$ws.Range("D17").Value = 135
$ws.Range("E17").Value = -12.592592592592
$ws.Range("F17").Value = 526
$ws.Range("H17").Value = 4.990019960079
$ws.Range("I17").Value = 6812
$ws.Range("J17").Value = 5859
$ws.Range("K17").Value = 16.26557433009
$ws.Range("L17").Value = 24.739058780443
$ws.Range("M17").Value = 63.161676646706
$ws.Range("N17").Value = -19.36553030303
$ws.Range("C18").Value = 38
$ws.Range("D18").Value = 57
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 191
$ws.Range("G18").Value = 206
$ws.Range("H18").Value = -7.281553398058
$ws.Range("I18").Value = 2703
$ws.Range("J18").Value = 2063
$ws.Range("K18").Value = 31.022782355792
$ws.Range("L18").Value = 7.818109293976
$ws.Range("M18").Value = -13.393143223325
$ws.Range("N18").Value = -84.635934746774
$ws.Range("C19").Value = 144
$ws.Range("D19").Value = 169
$ws.Range("E19").Value = -14.792899408284
$ws.Range("F19").Value = 605
$ws.Range("G19").Value = 668
$ws.Range("H19").Value = -9.43113772455
$ws.Range("I19").Value = 7454
$ws.Range("J19").Value = 6343
$ws.Range("K19").Value = 17.515371275421
$ws.Range("L19").Value = 31.255502729353
$ws.Range("M19").Value = 70.806599450045
$ws.Range("N19").Value = 6.516147470705
$ws.Range("C20").Value = 91
$ws.Range("D20").Value = 85
$ws.Range("E20").Value = 7.058823529411
$ws.Range("F20").Value = 315
$ws.Range("G20").Value = 331
$ws.Range("H20").Value = -4.833836858006
$ws.Range("I20").Value = 3616
$ws.Range("J20").Value = 2818
$ws.Range("K20").Value = 28.317955997161
$ws.Range("L20").Value = 89.220303506017
$ws.Range("M20").Value = 84.583971413986
$ws.Range("N20").Value = -74.876676162023
$ws.Range("D21").Value = 569
$ws.Range("E21").Value = -15.114235500878
$ws.Range("F21").Value = 2065
$ws.Range("G21").Value = 2162
$ws.Range("H21").Value = -4.486586493987
$ws.Range("I21").Value = 25811
$ws.Range("J21").Value = 21241
$ws.Range("K21").Value = 21.514994585942
$ws.Range("L21").Value = 34.481321315062
$ws.Range("M21").Value = 40.997487162678
$ws.Range("N21").Value = -59.539447901807
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = -20
$ws.Range("F22").Value = 22
$ws.Range("G22").Value = 47
$ws.Range("H22").Value = -53.191489361702
$ws.Range("I22").Value = 336
$ws.Range("J22").Value = 267
$ws.Range("K22").Value = 25.842696629213
$ws.Range("L22").Value = 10.89108910891
$ws.Range("M22").Value = 15.068493150684
$ws.Range("C23").Value = 25
$ws.Range("D23").Value = 29
$ws.Range("E23").Value = -13.793103448275
$ws.Range("F23").Value = 115
$ws.Range("G23").Value = 112
$ws.Range("H23").Value = 2.678571428571
$ws.Range("I23").Value = 1471
$ws.Range("J23").Value = 1171
$ws.Range("K23").Value = 25.619128949615
$ws.Range("L23").Value = 30.871886120996
$ws.Range("M23").Value = 46.222664015904
$ws.Range("C24").Value = 373
$ws.Range("D24").Value = 367
$ws.Range("E24").Value = 1.634877384196
$ws.Range("F24").Value = 1346
$ws.Range("G24").Value = 1265
$ws.Range("H24").Value = 6.403162055335
$ws.Range("I24").Value = 17168
$ws.Range("J24").Value = 12453
$ws.Range("K24").Value = 37.862362482935
$ws.Range("L24").Value = 30.90354555852
$ws.Range("M24").Value = 42.864275609553
$ws.Range("C25").Value = 159
$ws.Range("D25").Value = 193
$ws.Range("E25").Value = -17.61658031088
$ws.Range("F25").Value = 681
$ws.Range("G25").Value = 750
$ws.Range("H25").Value = -9.2
$ws.Range("I25").Value = 9206
$ws.Range("J25").Value = 8173
$ws.Range("K25").Value = 12.639177780496
$ws.Range("L25").Value = 11.506782945736
$ws.Range("M25").Value = -11.241804859236
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -26.666666666666
$ws.Range("F26").Value = 34
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 602
$ws.Range("J26").Value = 567
$ws.Range("K26").Value = 6.172839506172
$ws.Range("L26").Value = 24.123711340206
$ws.Range("C27").Value = 15
$ws.Range("D27").Value = 19
$ws.Range("E27").Value = -21.052631578947
$ws.Range("F27").Value = 56
$ws.Range("G27").Value = 77
$ws.Range("H27").Value = -27.272727272727
$ws.Range("I27").Value = 850
$ws.Range("J27").Value = 868
$ws.Range("K27").Value = -2.073732718894
$ws.Range("L27").Value = 29.573170731707
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 11
$ws.Range("E28").Value = -45.454545454545
$ws.Range("I28").Value = 447
$ws.Range("J28").Value = 570
$ws.Range("K28").Value = -21.578947368421
$ws.Range("L28").Value = 1.822323462414
$ws.Range("M28").Value = -1.541850220264
$ws.Range("N28").Value = -66.691505216095
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 10
$ws.Range("E29").Value = -40
$ws.Range("G29").Value = 44
$ws.Range("H29").Value = -36.363636363636
$ws.Range("I29").Value = 382
$ws.Range("J29").Value = 481
$ws.Range("K29").Value = -20.58212058212
$ws.Range("L29").Value = 5.524861878453
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -68.455821635012
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = -66.666666666666
$ws.Range("J30").Value = 46
$ws.Range("K30").Value = -8.695652173913
